$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Controller12927"
$ws.Range("B2").Value = "Automation12927"
$ws.Range("C2").Value = "controllerautomation12927@gmail.com"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "12927"
$ws.Range("D2").Style = "Normal"

$ws.Range("A3").Value = "Controller40980"
$ws.Range("B3").Value = "Automation40980"
$ws.Range("C3").Value = "controllerautomation40980@gmail.com"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40980"
$ws.Range("D3").Style = "Normal"

$ws.Range("A4").Value = "Controller42031"
$ws.Range("B4").Value = "Automation42031"
$ws.Range("C4").Value = "controllerautomation42031@gmail.com"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "42031"
$ws.Range("D4").Style = "Normal"

$ws.Range("A5").Value = "Controller00433"
$ws.Range("B5").Value = "Automation00433"
$ws.Range("C5").Value = "controllerautomation00433@gmail.com"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "00433"
$ws.Range("D5").Style = "Normal"

$ws.Range("A6").Value = "Controller74772"
$ws.Range("B6").Value = "Automation74772"
$ws.Range("C6").Value = "controllerautomation74772@gmail.com"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "74772"
$ws.Range("D6").Style = "Normal"

$ws.Range("A7").Value = "Controller32431"
$ws.Range("B7").Value = "Automation32431"
$ws.Range("C7").Value = "controllerautomation32431@gmail.com"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "32431"
$ws.Range("D7").Style = "Normal"

$ws.Range("A8").Value = "Controller59839"
$ws.Range("B8").Value = "Automation59839"
$ws.Range("C8").Value = "controllerautomation59839@gmail.com"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "59839"
$ws.Range("D8").Style = "Normal"

$ws.Range("A9").Value = "Controller57802"
$ws.Range("B9").Value = "Automation57802"
$ws.Range("C9").Value = "controllerautomation57802@gmail.com"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "57802"
$ws.Range("D9").Style = "Normal"

$ws.Range("A10").Value = "Controller62324"
$ws.Range("B10").Value = "Automation62324"
$ws.Range("C10").Value = "controllerautomation62324@gmail.com"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "62324"
$ws.Range("D10").Style = "Normal"

$ws.Range("A11").Value = "Controller78047"
$ws.Range("B11").Value = "Automation78047"
$ws.Range("C11").Value = "controllerautomation78047@gmail.com"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "78047"
$ws.Range("D11").Style = "Normal"

$ws.Range("A12").Value = "Controller97459"
$ws.Range("B12").Value = "Automation97459"
$ws.Range("C12").Value = "controllerautomation97459@gmail.com"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97459"
$ws.Range("D12").Style = "Normal"

$ws.Range("A13").Value = "Controller51803"
$ws.Range("B13").Value = "Automation51803"
$ws.Range("C13").Value = "controllerautomation51803@gmail.com"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "51803"
$ws.Range("D13").Style = "Normal"

$ws.Range("A14").Value = "Controller67974"
$ws.Range("B14").Value = "Automation67974"
$ws.Range("C14").Value = "controllerautomation67974@gmail.com"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "67974"
$ws.Range("D14").Style = "Normal"

$ws.Range("A15").Value = "Controller87666"
$ws.Range("B15").Value = "Automation87666"
$ws.Range("C15").Value = "controllerautomation87666@gmail.com"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "87666"
$ws.Range("D15").Style = "Normal"

$ws.Range("A16").Value = "Controller16480"
$ws.Range("B16").Value = "Automation16480"
$ws.Range("C16").Value = "controllerautomation16480@gmail.com"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16480"
$ws.Range("D16").Style = "Normal"

$ws.Range("A17").Value = "Controller12533"
$ws.Range("B17").Value = "Automation12533"
$ws.Range("C17").Value = "controllerautomation12533@gmail.com"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "12533"
$ws.Range("D17").Style = "Normal"

$ws.Range("A18").Value = "Controller23544"
$ws.Range("B18").Value = "Automation23544"
$ws.Range("C18").Value = "controllerautomation23544@gmail.com"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "23544"
$ws.Range("D18").Style = "Normal"

$ws.Range("A19").Value = "Controller36449"
$ws.Range("B19").Value = "Automation36449"
$ws.Range("C19").Value = "controllerautomation36449@gmail.com"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36449"
$ws.Range("D19").Style = "Normal"

$ws.Range("A20").Value = "Controller94654"
$ws.Range("B20").Value = "Automation94654"
$ws.Range("C20").Value = "controllerautomation94654@gmail.com"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "94654"
$ws.Range("D20").Style = "Normal"

$ws.Range("A21").Value = "Controller33228"
$ws.Range("B21").Value = "Automation33228"
$ws.Range("C21").Value = "controllerautomation33228@gmail.com"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "33228"
$ws.Range("D21").Style = "Normal"

$ws.Range("A22").Value = "Controller16378"
$ws.Range("B22").Value = "Automation16378"
$ws.Range("C22").Value = "controllerautomation16378@gmail.com"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16378"
$ws.Range("D22").Style = "Normal"

$ws.Range("A23").Value = "Controller92436"
$ws.Range("B23").Value = "Automation92436"
$ws.Range("C23").Value = "controllerautomation92436@gmail.com"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "92436"
$ws.Range("D23").Style = "Normal"

$ws.Range("A24").Value = "Controller94809"
$ws.Range("B24").Value = "Automation94809"
$ws.Range("C24").Value = "controllerautomation94809@gmail.com"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "94809"
$ws.Range("D24").Style = "Normal"

$ws.Range("A25").Value = "Controller10331"
$ws.Range("B25").Value = "Automation10331"
$ws.Range("C25").Value = "controllerautomation10331@gmail.com"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10331"
$ws.Range("D25").Style = "Normal"

$ws.Range("A26").Value = "Controller03781"
$ws.Range("B26").Value = "Automation03781"
$ws.Range("C26").Value = "controllerautomation03781@gmail.com"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "03781"
$ws.Range("D26").Style = "Normal"

$ws.Range("A27").Value = "Controller08936"
$ws.Range("B27").Value = "Automation08936"
$ws.Range("C27").Value = "controllerautomation08936@gmail.com"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "08936"
$ws.Range("D27").Style = "Normal"

$ws.Range("A28").Value = "Controller39298"
$ws.Range("B28").Value = "Automation39298"
$ws.Range("C28").Value = "controllerautomation39298@gmail.com"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39298"
$ws.Range("D28").Style = "Normal"

$ws.Range("A29").Value = "Controller92304"
$ws.Range("B29").Value = "Automation92304"
$ws.Range("C29").Value = "controllerautomation92304@gmail.com"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "92304"
$ws.Range("D29").Style = "Normal"

$ws.Range("A30").Value = "Controller01737"
$ws.Range("B30").Value = "Automation01737"
$ws.Range("C30").Value = "controllerautomation01737@gmail.com"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "01737"
$ws.Range("D30").Style = "Normal"

$ws.Range("A31").Value = "Controller89338"
$ws.Range("B31").Value = "Automation89338"
$ws.Range("C31").Value = "controllerautomation89338@gmail.com"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "89338"
$ws.Range("D31").Style = "Normal"

